# Leave Card update — 10/3/2023 3:18 PM
# Adds six more monthly VL-earned rows (Mar-Aug 2023), plus a VL(10-0-0)
# leave entry in September 2023, and fills in the remaining month-start
# dates for Oct/Nov 2023.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Rows 61-66: month-start date in column A, 1.25 EARNED in column C ---
$ws.Range("A61").Value = 44986   # 3/1/2023
$ws.Range("C61").Value = 1.25

$ws.Range("A62").Value = 45017   # 4/1/2023
$ws.Range("C62").Value = 1.25

$ws.Range("A63").Value = 45047   # 5/1/2023
$ws.Range("C63").Value = 1.25

$ws.Range("A64").Value = 45078   # 6/1/2023
$ws.Range("C64").Value = 1.25

$ws.Range("A65").Value = 45108   # 7/1/2023
$ws.Range("C65").Value = 1.25

$ws.Range("A66").Value = 45139   # 8/1/2023
$ws.Range("C66").Value = 1.25

# --- Row 67: a VL(10-0-0) leave taken 10/18-31/2023, charged 10 days ---
$ws.Range("A67").Value = 45170   # 9/1/2023
$ws.Range("B67").Value = "VL(10-0-0)"
$ws.Range("D67").Value = 10
$ws.Range("K67").Value = "10/18-31/2023"

# --- Rows 68-69: month-start date only ---
$ws.Range("A68").Value = 45200   # 10/1/2023
$ws.Range("A69").Value = 45231   # 11/1/2023

# Reflect the author's last on-screen selection (bottom pane, cell E67).
[void]$ws.Activate()
[void]$ws.Range("E67").Select()
